# The commit resizes the "Directory / Description" table on the
# "Filesystem Hierarchy Standard (FHS)" slide (sldId 490, slide 22 in
# the current slide order): the first ("Directory") column was narrowed
# by dragging its left border inward, which both shrinks the column and
# the overall table width while shifting the table's left edge right by
# the same amount; the second ("Description") column keeps its width.
#
#   off  x: 1156008 -> 1159497   (+3489 EMU)
#   ext  cx: 7313343 -> 7309854   (-3489 EMU)
#   gridCol0 w: 1932879 -> 1929390   (-3489 EMU)
#   gridCol1 w: 5380464 -> 5380464   (unchanged)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)

# Locate the table graphic frame (shape id 5, "표 4") on the slide.
$tblShape = $s.Shapes.Item(3)
$tbl = $tblShape.Table

# Narrow the first column by ~0.27466pt (3489 EMU / 12700); PowerPoint
# keeps the table's right edge fixed and recomputes the overall table
# width from the column widths, so this alone drives ext/cx and
# gridCol[0] to their target values while leaving gridCol[1] untouched.
$tbl.Columns.Item(1).Width = 151.92047244094488

# Dragging the left border also nudges the frame's left edge right by
# the same amount that the first column shrank.
$tblShape.Left = 91.29901885986328
